# Cryptos list refresh -- Mon Jul  1 22:56:22 UTC 2024 (GitHub Actions bot run).
# Re-writes the "Price" (column D) and "Volume(1h)" (column E) text values on the
# active sheet to match the latest coinranking.com snapshot. Both columns hold
# plain text in the source data (e.g. "3.438.73" uses '.' as a thousands marker,
# and the volume column keeps its padding spaces), so for Price values that would
# otherwise be auto-detected as a number by Excel we briefly mark the cell as Text,
# write the literal string, then restore the cell's original (Normal) style so no
# formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.937.53'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '3.438.73'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.59'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').Value = '3.438.15'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.78'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.53%  '
$ws.Range('E11').Value = '  -1.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.404'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.17%  '
$ws.Range('D13').Value = '4.027.23'
$ws.Range('E13').Value = '  -0.22%  '
$ws.Range('E14').Value = '  +2.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.93'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.91%  '
$ws.Range('D16').Value = '3.435.36'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('D18').Value = '62.970.25'
$ws.Range('E18').Value = '  +0.32%  '
$ws.Range('E19').Value = '  +2.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.39'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('E21').Value = '  -1.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '385.53'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.66%  '
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.45'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.54%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').Value = '3.590.51'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('E27').Value = '  -3.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.57'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.50%  '
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.09'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.10%  '
$ws.Range('E32').Value = '  -1.71%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.25'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.44%  '
$ws.Range('E35').Value = '  -8.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.28'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.61%  '
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.59'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '31.53'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '168.91'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('D41').Value = '3.475.30'
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0765'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.788'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.37'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.40%  '
$ws.Range('E45').Value = '  -0.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.20'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.35'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.04%  '
$ws.Range('D48').Value = '2.570.57'
$ws.Range('E48').Value = '  +1.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.27'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.79'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.69'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.29%  '
